# Atualização de bases das ligas, do dia: 20-06-2024 às 20:11
#
# The source data rows got re-sorted/re-paired (adjacent matches on the same
# date swapped places, and one 5-row block rotated by one position). Column
# A (row index) and the already-identical C (Div) / D (Date) columns stay
# put; every other column (B id, E HomeTeam, F AwayTeam, G..AD stats/odds)
# moves with its match record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that travel with each match record (everything except A/C/D).
$cols = @("B","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")

# Each inner array is a cyclic group of worksheet row numbers: the data that
# was in row[i] moves into row[i-1] (i.e. row[i] receives what used to be in
# row[i+1], wrapping around).
$groups = @(
    @(116, 117),
    @(155, 156),
    @(190, 191),
    @(200, 201),
    @(208, 209),
    @(262, 263),
    @(302, 303, 304, 305, 306)
)

foreach ($group in $groups) {
    $n = $group.Length

    # Snapshot every cell we might touch before writing any of them back,
    # since rows within a group overwrite each other.
    $snapshot = @()
    foreach ($r in $group) {
        $rowVals = @{}
        foreach ($c in $cols) {
            $rowVals[$c] = $ws.Range("$c$r").Value()
        }
        $snapshot += ,$rowVals
    }

    for ($i = 0; $i -lt $n; $i++) {
        $destRow = $group[$i]
        $srcVals = $snapshot[($i + 1) % $n]
        foreach ($c in $cols) {
            $ws.Range("$c$destRow").Value = $srcVals[$c]
        }
    }
}
